# Generate Report for Handoff
#
# File 86eabc96-df8b-43f1-a388-eeccc15753de.md has finished translation and
# is now ready to be handed off again. Update its status on every sheet and
# refresh the handoff/handback timestamps for the files whose report rows
# changed during this run (the file that just became ready, plus the file
# that was already ready and the file whose handback transform failed).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date -------
# Row 6  -> c2f9c434-7e2f-4ca2-9c29-2f130809b89a.md (Handback transform failed)
$overview.Range("D6").Value = "2016-18-12 00:18:15"

# Row 9  -> 86eabc96-df8b-43f1-a388-eeccc15753de.md (In Translation -> Ready for handoff)
$overview.Range("B9").Value = "Ready for handoff"
$overview.Range("C9").Value = "Ready for handoff"
$overview.Range("D9").Value = "2016-18-12 00:18:15"

# Row 10 -> c8018f2e-b428-4c51-b373-9ec6c5ca8a41.md (already Ready for handoff)
$overview.Range("D10").Value = "2016-18-12 00:18:15"

# --- zh-cn sheet: Status (C) + Latest Handback DateTime (E) ---------------
$zhcn.Range("E6").Value = "2016-03-12 00:18:12"

$zhcn.Range("C9").Value = "Ready for handoff"
$zhcn.Range("E9").Value = "2016-03-12 00:18:12"

$zhcn.Range("E10").Value = "2016-03-12 00:18:12"

# --- de-de sheet: Status (C) + Latest Handback DateTime (E) ---------------
$dede.Range("E6").Value = "2016-03-12 00:18:15"

$dede.Range("C9").Value = "Ready for handoff"
$dede.Range("E9").Value = "2016-03-12 00:18:15"

$dede.Range("E10").Value = "2016-03-12 00:18:15"
